$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @('21CRB01268', 'Bunner', 'POSSESSION DRUG PARAPHERNALIA', '2925.14(C)', 'M4', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21CRB01268', 'Bunner', 'No Operator License - Never Held', '4510.12(C)(1)', 'Unclassified Misdemeanor', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21CRB01268', 'Bunner', 'POSSESSION DRUG PARAPHERNALIA', '2925.14(C)', 'M4', 'Guilty', 'Guilty', '$ 0', '$ 0'),
    @('21CRB01268', 'Bunner', 'Driving Under Financial Responsibility Law Suspension', '4510.16', 'Unclassified Misdemeanor', 'Guilty', 'Guilty', '$ 0', '$ 0'),
    @('21CRB01268', 'Bunner', 'POSSESSION DRUG PARAPHERNALIA', '2925.14(C)', 'M4', 'Guilty', 'Guilty', '$ 0', '$ 0'),
    @('21CRB01268', 'Bunner', 'Driving Under Financial Responsibility Law Suspension', '4510.16', 'Unclassified Misdemeanor', 'Guilty', 'Guilty', '$ 0', '$ 0'),
    @('21CRB01268', 'Bunner', 'POSSESSION DRUG PARAPHERNALIA', '2925.14(C)', 'M4', 'Guilty', 'Guilty', '$ 0', '$ 0'),
    @('21CRB01268', 'Bunner', 'Driving Under Financial Responsibility Law Suspension', '4510.16', 'Unclassified Misdemeanor', 'Guilty', 'Guilty', '$ 0', '$ 0'),
    @('21CRB01268', 'Bunner', 'POSSESSION DRUG PARAPHERNALIA', '2925.14(C)', 'M4', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21CRB01291', 'Hemmeter', "PERMISSION REQ'D TO USE LICENSED DOCK", '1501:46-12-04', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21CRB01291', 'Hemmeter', 'No Operator License - Expired', '4510.12(C)(2)', 'Minor Misdemeanor', 'No Contest', 'Guilty', '$ 0', '$ 0')
)

$startRow = 550
$endRow = $startRow + $rows.Count - 1

# Force the target range to Text format so numeric-looking strings
# (e.g. "4510.16", "$ 0") are stored as text, matching the source data.
$rangeAddr = "A" + $startRow + ":I" + $endRow
$ws.Range($rangeAddr).NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}
